$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (C2:K2) to 1
$ws.Range("C2:K2").Value = 1

# Update the active selection to K3
$ws.Range("K3").Select()
